$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old "Docentes responsaveis" names (rows 13 & 14) are removed; the
# sheet's later content shifts up by two rows to take their place, and new
# rows (resumido/short syllabus/programa/syllabus/avaliacao/metodo) are
# inserted in what used to be rows 13-14's slot as part of that shift.
$ws.Rows.Item(13).Delete()
$ws.Rows.Item(13).Delete()

# --- Text updates -----------------------------------------------------
# Plain text replacements (not number/date-like, safe to assign directly).
$ws.Range("B10").Value = "471420 - Carlos Antonio Reis Pereira Baptista"
$ws.Range("C10").Value = "471420 - Carlos Antonio Reis Pereira Baptista"

$ws.Range("B15").Value = "471420 - Carlos Antonio Reis Pereira Baptista"
$ws.Range("C15").Value = "471420 - Carlos Antonio Reis Pereira Baptista"

$ws.Range("B18").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C18").Value = "519033 - Carlos Yujiro Shigue"

$ws.Range("B19").Value = "Este curso deverá conter avaliações escritas e desenvolvimento de Estudo de Casos ou Projetos na área de Engenharia de Materiais. Sendo necessário aplicar pelo menos dois tipos de avaliações diferentes."
$ws.Range("C19").Value = "Este curso deverá conter avaliações escritas e desenvolvimento de Estudo de Casos ou Projetos na área de Engenharia de Materiais. Sendo necessário aplicar pelo menos dois tipos de avaliações diferentes."

$ws.Range("B20").Value = "A média do semestre será computada com base na relação: M=(A1+A2)/2"
$ws.Range("C20").Value = "A média do semestre será computada com base na relação: M=(A1+A2)/2"

$ws.Range("B21").Value = "Não cabe recuperação."
$ws.Range("C21").Value = "Não cabe recuperação."

# "01/01/2022" is date-like text (stored as a shared string, not a real
# date, in the source file) -- a direct .Value assignment of that literal
# would be auto-converted to a date serial by Excel's type inference, and
# picking up a number format would also fabricate a new (unused) style
# entry. Copy the text, verbatim, from the cell that already holds it
# (B8) via PasteSpecial(xlPasteValues = -4163), which carries over the
# original text type without touching the destination cell's own style.
$ws.Range("B8").Copy()
$ws.Range("B13").PasteSpecial(-4163)
$ws.Range("B8").Copy()
$ws.Range("C13").PasteSpecial(-4163)
$excel.CutCopyMode = $false
